# Applies the row permutation/content rotation described in the commit diff.
# Rows 4-13 on the active sheet ("Artfynd") have their data shifted: each
# row's record (taxon info, coordinates, growth-substrate note, etc.) is
# replaced by another row's original record, per the mapping:
#   new4<-old5, new5<-old13, new6<-old4, new7<-old6, new8<-old7,
#   new9<-old8, new10<-old9, new11<-old10, new12<-old11, new13<-old12
# Columns C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY
# are identical across all these rows, so they require no edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # -- Row 4 --
    $ws.Range("A4").Value = 111473782
    $ws.Range("B4").Value = 89183
    $ws.Range("D4").Value = "LC"
    $ws.Range("E4").Value = 3215
    $ws.Range("F4").Value = "Rödgul trumpetsvamp"
    $ws.Range("G4").Value = "Craterellus lutescens"
    $ws.Range("H4").Value = "(Fr.) Fr."
    $ws.Range("Q4").Value = 704171.5165585374
    $ws.Range("R4").Value = 6572850.843097115
    $ws.Range("AO4").Value = ""

    # -- Row 5 --
    $ws.Range("A5").Value = 111473779
    $ws.Range("B5").Value = 89425
    $ws.Range("D5").Value = "NT"
    $ws.Range("E5").Value = 5442
    $ws.Range("F5").Value = "Tallticka"
    $ws.Range("G5").Value = "Porodaedalea pini"
    $ws.Range("H5").Value = "(Brot.) Murrill"
    $ws.Range("Q5").Value = 704193.4830821306
    $ws.Range("R5").Value = 6572948.378178579
    $ws.Range("AO5").Value = "gammeltall"

    # -- Row 6 --
    $ws.Range("A6").Value = 111473775
    $ws.Range("B6").Value = 89405
    $ws.Range("D6").Value = "NT"
    $ws.Range("E6").Value = 1202
    $ws.Range("F6").Value = "Ullticka"
    $ws.Range("G6").Value = "Phellinidium ferrugineofuscum"
    $ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
    $ws.Range("Q6").Value = 703969.3444121893
    $ws.Range("R6").Value = 6572791.287347207

    # -- Row 7 --
    $ws.Range("A7").Value = 111473783
    $ws.Range("B7").Value = 89686
    $ws.Range("D7").Value = "NT"
    $ws.Range("E7").Value = 658
    $ws.Range("F7").Value = "Rosenticka"
    $ws.Range("G7").Value = "Rhodofomes roseus"
    $ws.Range("H7").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
    $ws.Range("Q7").Value = 703998.3853129407
    $ws.Range("R7").Value = 6572852.813158008

    # -- Row 8 --
    $ws.Range("A8").Value = 111473774
    $ws.Range("B8").Value = 89405
    $ws.Range("D8").Value = "NT"
    $ws.Range("E8").Value = 1202
    $ws.Range("F8").Value = "Ullticka"
    $ws.Range("G8").Value = "Phellinidium ferrugineofuscum"
    $ws.Range("H8").Value = "(P.Karst.) Fiasson & Niemelä"
    $ws.Range("Q8").Value = 703999.5190368021
    $ws.Range("R8").Value = 6572850.823973293
    $ws.Range("M8").Value = ""
    $ws.Range("AO8").Value = "granlåga"

    # -- Row 9 --
    $ws.Range("A9").Value = 111473792
    $ws.Range("B9").Value = 5113
    $ws.Range("D9").Value = "LC"
    $ws.Range("E9").Value = 100526
    $ws.Range("F9").Value = "Bronshjon"
    $ws.Range("G9").Value = "Callidium coriaceum"
    $ws.Range("H9").Value = "Paykull, 1800"
    $ws.Range("Q9").Value = 703965.55072247
    $ws.Range("R9").Value = 6572785.445717536
    $ws.Range("M9").Value = "äldre gnagspår"
    $ws.Range("AO9").Value = "torrgran"

    # -- Row 10 --
    $ws.Range("A10").Value = 111473773
    $ws.Range("B10").Value = 89405
    $ws.Range("D10").Value = "NT"
    $ws.Range("E10").Value = 1202
    $ws.Range("F10").Value = "Ullticka"
    $ws.Range("G10").Value = "Phellinidium ferrugineofuscum"
    $ws.Range("H10").Value = "(P.Karst.) Fiasson & Niemelä"
    $ws.Range("Q10").Value = 704016.0051346947
    $ws.Range("R10").Value = 6572801.994589122
    $ws.Range("AO10").Value = "granlåga"

    # -- Row 11 --
    $ws.Range("A11").Value = 111473791
    $ws.Range("B11").Value = 93289
    $ws.Range("D11").Value = "LC"
    $ws.Range("E11").Value = 2170
    $ws.Range("F11").Value = "Flagellkvastmossa"
    $ws.Range("G11").Value = "Dicranum flagellare"
    $ws.Range("H11").Value = "Hedw."
    $ws.Range("Q11").Value = 704004.9502936595
    $ws.Range("R11").Value = 6572835.740028554
    $ws.Range("AO11").Value = "låga av tall"

    # -- Row 12 --
    $ws.Range("A12").Value = 111473793
    $ws.Range("B12").Value = 93388
    $ws.Range("D12").Value = "LC"
    $ws.Range("E12").Value = 2180
    $ws.Range("F12").Value = "Blåmossa"
    $ws.Range("G12").Value = "Leucobryum glaucum"
    $ws.Range("H12").Value = "(Hedw.) Ångstr."
    $ws.Range("Q12").Value = 703959.3331032015
    $ws.Range("R12").Value = 6572805.612961343
    $ws.Range("AO12").Value = ""

    # -- Row 13 --
    $ws.Range("A13").Value = 111473776
    $ws.Range("B13").Value = 89405
    $ws.Range("D13").Value = "NT"
    $ws.Range("E13").Value = 1202
    $ws.Range("F13").Value = "Ullticka"
    $ws.Range("G13").Value = "Phellinidium ferrugineofuscum"
    $ws.Range("H13").Value = "(P.Karst.) Fiasson & Niemelä"
    $ws.Range("Q13").Value = 703970.8884549731
    $ws.Range("R13").Value = 6572810.333898042
    $ws.Range("AO13").Value = "granlåga"

Write-Output "Rows 4-13 updated"
